# Update crypto price/volume figures and reorder Stacks/Fetch.AI/Bittensor rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.388.02"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "3.599.77"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "200.76"
$ws.Range("E5").Value = "  +2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "593.79"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +6.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.645"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.47"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000301"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "693.61"
$ws.Range("E14").Value = "  +16.63%  "
$ws.Range("D15").Value = "4.168.91"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").Value = "70.489.00"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.75"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.07"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "3.596.63"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.37"
$ws.Range("E22").Value = "  +3.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "110.64"
$ws.Range("E23").Value = "  +7.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.29"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.54"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.53"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.02"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.04"
$ws.Range("E29").Value = "  +5.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.61"
$ws.Range("E30").Value = "  +4.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.46"
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.07"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.28"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.68"
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("D36").Value = "0.0₃0849"
$ws.Range("E36").Value = "  +3.35%  "
$ws.Range("D37").Value = "3.827.28"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.63"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.02"
$ws.Range("E40").Value = "  -5.55%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "510.09"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.69"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.382"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.138"
$ws.Range("E44").Value = "  +3.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0472"
$ws.Range("E45").Value = "  +4.84%  "
$ws.Range("E46").Value = "  +8.57%  "
$ws.Range("E47").Value = "  +5.27%  "
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.66"
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.81"
$ws.Range("E51").Value = "  +19.82%  "
